$d = $word.ActiveDocument

# Smart (curly) right single quotation mark used in "It's" -> "It’s"
$rsquo = [char]0x2019

# The document currently has a single paragraph: "I love working with git."
$firstPara = $d.Paragraphs.First

# Insert a new paragraph BEFORE the existing one and fill it in.
$beforeRange = $firstPara.Range
$beforeRange.InsertParagraphBefore()
$d.Paragraphs.First.Range.Text = "I love working with GIT. It${rsquo}s simply fast and superb version control system."

# Insert a new paragraph AFTER the existing (still-middle) one and fill it in.
$lastPara = $d.Paragraphs.Last
$afterRange = $lastPara.Range
$afterRange.InsertParagraphAfter()
$d.Paragraphs.Last.Range.Text = "It increases my productivity on multiple folds when working with files which has frequency of changes."
